$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "c1"
$ws.Range("A3").Value = "c2"
$ws.Range("A4").Value = "c3"
$ws.Range("B2").Value = "s1"
$ws.Range("B3").Value = "s2"
$ws.Range("B4").Value = "s3"

$ws.Range("B5").Select()
